# "adding averages and more checks"
#
# Style change: the bold title font (size 14, black) and the bold header
# font used on the blue header band no longer differ by size - the title
# font loses its explicit 14pt size and both the title style and the
# header style now share one bold white font.
#
# Data change (Training Dashboard, row 3): the "PERIOD TO EXPIRE" and
# "LAST UPDATE" values were recalculated against a later check date.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Title cell (row 1) - bold, white, default size (was bold/size14/black)
    $titleRange = $ws.Range("A1")
    $titleRange.Font.Bold = $true
    $titleRange.Font.Size = 11
    $titleRange.Font.Color = 16777215

    # Header band (row 2) - bold, white (was bold/black) on the blue fill
    $headerRange = $ws.Rows.Item(2)
    $headerRange.Font.Bold = $true
    $headerRange.Font.Color = 16777215
}

# Training Dashboard sheet - refresh the expiry check for row 3
$trainingSheet = $wb.Worksheets.Item("Training Dashboard")
$trainingSheet.Range("H3").Value = -55

# Write I3 as literal text ("16-Sep-2025") rather than letting Excel parse
# it as a date serial: format as Text first, enter the value, then
# re-apply H3's (unchanged) cell format so I3 keeps its original look.
$trainingSheet.Range("I3").NumberFormat = "@"
$trainingSheet.Range("I3").Value = "16-Sep-2025"
$trainingSheet.Range("H3").Copy()
$trainingSheet.Range("I3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
